$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A27").Value = "Mayflies"
$ws.Range("B27").Value = "Andrew O'Hagan"
$ws.Range("C27").Value = "Literature; Scottish"
$ws.Range("D27").Value = 8.99

$ws.Range("A28").Value = "Later"
$ws.Range("B28").Value = "Stephen King"
$ws.Range("C28").Value = "Horror; Thriller"
$ws.Range("D28").Value = 8.99

$ws.Range("I25").Select()
